$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 2, pushing existing rows 2-3 down to 4-5
$ws.Range("A2:A3").EntireRow.Insert()

# Insert copies the formatting of the row above into the new rows. Column A
# should keep the bold/boxed header-style (matching A4/A5 below), while the
# data columns B:P should have no special formatting at all.
$ws.Range("A4:A5").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$ws.Range("B2:P3").ClearFormats()

# New row 2: 新视云 (D2 stays empty, as in the original data)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "新视云"
$ws.Range("C2").Value = "雨花台"
$ws.Range("E2").Value = "Java"
$ws.Range("F2").Value = "9:00-17:30"
$ws.Range("G2").Value = "1h"
$ws.Range("H2").Value = "看部门，业务部门偶尔加班，技术支持部门基本不加班"
$ws.Range("I2").Value = "基数5k，比例8%"
$ws.Range("J2").Value = "固定13薪"
$ws.Range("K2").Value = "3年合同，试用期总共6个月，前三个月8折，后三个月全薪"
$ws.Range("L2").Value = "配笔记本+显示器"
$ws.Range("M2").Value = "5天年假+5天带薪病假（入职自动折算当年年假）"
$ws.Range("N2").Value = "不打卡"

# New row 3: 华为 (D3 stays empty, as in the original data)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "华为"
$ws.Range("C3").Value = "华为南研所"
$ws.Range("E3").Value = "Java"
$ws.Range("F3").Value = "9:00"
$ws.Range("G3").Value = "12:00-13:40"
$ws.Range("H3").Value = "看部门情况。好部门：124加班8：30，35正常下班,差部门：天天11点以后"
$ws.Range("I3").Value = "基础工资的5%"
$ws.Range("J3").Value = "看部门盈利情况和个人绩效定"
$ws.Range("K3").Value = "试用期6个月，100%工资不打折"
$ws.Range("L3").Value = "配win台式机+双屏"
$ws.Range("M3").Value = "没签奋斗协议的5天，但一般不给休，第二年可以换成钱。签了的自愿放弃年假了"
$ws.Range("N3").Value = "必须按时打卡"

# Update the index column for the shifted-down rows (now 4 and 5)
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
